# Marksheet update: student answers were recorded (the student was no longer
# "Absent"), so the summary counts/scores and the per-question "Student Ans"
# columns need to be filled in. The sheet also tracked a 3rd question column
# (G:H) that is no longer used and is removed, along with the now-unused
# tail of the 2nd question column (D:E beyond row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block (rows 10-12) ----
# Row headers (No./Marking/Total) pick up the section title style.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# No. row: Right / Wrong / Not Attempt / Max
$ws.Range("B10").Value = 21
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 28

# Marking row: points per right / wrong answer
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Total row: marks scored / lost, and the "scored/max" summary
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "82/112"

# ---- Per-question "Student Ans" columns ----
# Column A holds the student's answer for question block 1 (rows 16-40);
# column B already holds the correct answer. Style reflects correctness:
# correctStyle (matches), incorrectStyle (mismatch) or normalStyle (blank).
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A17").Value = "Option D"
$ws.Range("A17").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A20").Value = "Option B"
$ws.Range("A20").Style = "correctStyle"
$ws.Range("A21").Value = "Option D"
$ws.Range("A21").Style = "incorrectStyle"
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A23").Style = "normalStyle"
$ws.Range("A24").Style = "normalStyle"
$ws.Range("A25").Value = "Option A"
$ws.Range("A25").Style = "correctStyle"
$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A28").Style = "normalStyle"
$ws.Range("A29").Value = "Option D"
$ws.Range("A29").Style = "correctStyle"
$ws.Range("A30").Value = "Option B"
$ws.Range("A30").Style = "correctStyle"
$ws.Range("A31").Value = "Option D"
$ws.Range("A31").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"
$ws.Range("A34").Style = "correctStyle"
$ws.Range("A35").Style = "normalStyle"
$ws.Range("A36").Value = "Option A"
$ws.Range("A36").Style = "correctStyle"
$ws.Range("A37").Value = "Option A"
$ws.Range("A37").Style = "correctStyle"
$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Style = "normalStyle"

# Column D holds the student's answer for question block 2, but only the
# first three questions of that block (rows 16-18) were attempted; the rest
# (rows 19-40, column D and E) are no longer used and get cleared below.
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value = "Option A"
$ws.Range("D17").Style = "incorrectStyle"
$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"

# ---- Drop unused columns/rows ----
# Question block 2 tail (no longer attempted) and question block 3 (G:H)
# are fully removed so the used range shrinks back to A:E.
$ws.Range("D19:E40").Clear()
$ws.Range("G15:H40").Clear()
